$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Row 43: I am C-3PO
$ws.Range("A43").Value = "I am C-3PO"
$ws.Range("B43").Value = "Anthony Daniels"
$ws.Range("C43").Value = 43912
$ws.Range("D43").Value = 43913
$ws.Range("E43").Value = "star wars;hollywood;movies;acting;c-3po"
$ws.Range("F43").Value = "Audio"
$ws.Range("G43").Value = "9 Hours 34 Mins"

# Row 44: Deep Learning
$ws.Range("A44").Value = "Deep Learning"
$ws.Range("B44").Value = "John D. Kellerher"
$ws.Range("C44").Value = 43911
$ws.Range("D44").Value = 43913
$ws.Range("E44").Value = "deep learning;machine learning;data science;neural networks"
$ws.Range("F44").Value = "Hard Copy"
$ws.Range("G44").Value = "250 Pages"

$ws.Range("C43:D44").NumberFormat = "m/d/yy"

$ws.Range("A45").Select()
$ws.Application.ActiveWindow.ScrollRow = 25
